$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column H (run cmds work) values for rows 2-9
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 0

# Update the selected cell/range in the sheet view
$ws.Range("H11").Select()
